$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header cells (A1:D1) to short snake_case column codes ---
# --- Title-case Spanish connector words (de/del/la/las/los/el/y) in state/municipality names ---
# --- Fix the one all-caps state name (GUANAJUATO -> Guanajuato) ---

$ws.Cells.Item(1, 1).Value = 'mx_state'
$ws.Cells.Item(1, 2).Value = 'mx_municipality'
$ws.Cells.Item(1, 3).Value = 'n_matriculas'
$ws.Cells.Item(1, 4).Value = 'pct_matriculas'
$ws.Cells.Item(17, 2).Value = 'Amatenango De La Frontera'
$ws.Cells.Item(19, 2).Value = 'Bejucal De Ocampo'
$ws.Cells.Item(21, 2).Value = 'Benemérito De Las Américas'
$ws.Cells.Item(24, 2).Value = 'Chiapa De Corzo'
$ws.Cells.Item(27, 2).Value = 'Comitán De Domínguez'
$ws.Cells.Item(43, 2).Value = 'Mazapa De Madero'
$ws.Cells.Item(48, 2).Value = 'Ocozocoautla De Espinosa'
$ws.Cells.Item(78, 2).Value = 'San Juan De Sabinas'
$ws.Cells.Item(84, 1).Value = 'Ciudad De México'
$ws.Cells.Item(88, 2).Value = 'Cuajimalpa De Morelos'
$ws.Cells.Item(103, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(106, 1).Value = 'Estado De México'
$ws.Cells.Item(106, 2).Value = 'Acambay De Ruíz Castañeda'
$ws.Cells.Item(108, 2).Value = 'Almoloya De Alquisiras'
$ws.Cells.Item(109, 2).Value = 'Almoloya Del Río'
$ws.Cells.Item(114, 2).Value = 'Atizapán De Zaragoza'
$ws.Cells.Item(119, 2).Value = 'Chapa De Mota'
$ws.Cells.Item(123, 2).Value = 'Coacalco De Berriozábal'
$ws.Cells.Item(130, 2).Value = 'Ecatepec De Morelos'
$ws.Cells.Item(134, 2).Value = 'Ixtapan De La Sal'
$ws.Cells.Item(145, 2).Value = 'Naucalpan De Juárez'
$ws.Cells.Item(154, 2).Value = 'San Felipe Del Progreso'
$ws.Cells.Item(155, 2).Value = 'San Martín De Las Pirámides'
$ws.Cells.Item(161, 2).Value = 'Tenango Del Valle'
$ws.Cells.Item(168, 2).Value = 'Tlalnepantla De Baz'
$ws.Cells.Item(172, 2).Value = 'Valle De Bravo'
$ws.Cells.Item(173, 2).Value = 'Valle De Chalco Solidaridad'
$ws.Cells.Item(174, 2).Value = 'Villa De Allende'
$ws.Cells.Item(175, 2).Value = 'Villa Del Carbón'
$ws.Cells.Item(181, 1).Value = 'Guanajuato'
$ws.Cells.Item(184, 2).Value = 'Apaseo El Alto'
$ws.Cells.Item(191, 2).Value = 'Dolores Hidalgo Cuna De La Independencia Nacional'
$ws.Cells.Item(195, 2).Value = 'Jaral Del Progreso'
$ws.Cells.Item(202, 2).Value = 'San Luis De La Paz'
$ws.Cells.Item(203, 2).Value = 'Santa Cruz De Juventino Rosas'
$ws.Cells.Item(206, 2).Value = 'Valle De Santiago'
$ws.Cells.Item(210, 2).Value = 'Acapulco De Juárez'
$ws.Cells.Item(212, 2).Value = 'Ajuchitlán Del Progreso'
$ws.Cells.Item(213, 2).Value = 'Alcozauca De Guerero'
$ws.Cells.Item(217, 2).Value = 'Atenango Del Río'
$ws.Cells.Item(218, 2).Value = 'Atlamajalcingo Del Monte'
$ws.Cells.Item(220, 2).Value = 'Atoyac De Álvarez'
$ws.Cells.Item(221, 2).Value = 'Ayutla De Los Libres'
$ws.Cells.Item(224, 2).Value = 'Chilapa De Álvarez'
$ws.Cells.Item(225, 2).Value = 'Chilpancingo De Los Bravo'
$ws.Cells.Item(226, 2).Value = 'Coahuayutla De José María Izazaga'
$ws.Cells.Item(229, 2).Value = 'Coyuca De Benítez'
$ws.Cells.Item(230, 2).Value = 'Coyuca De Catalán'
$ws.Cells.Item(234, 2).Value = 'Cuetzala Del Progreso'
$ws.Cells.Item(235, 2).Value = 'Cutzamala De Pinzón'
$ws.Cells.Item(239, 2).Value = 'Iguala De La Independencia'
$ws.Cells.Item(241, 2).Value = 'Zihuatanejo De Azueta'
$ws.Cells.Item(245, 2).Value = 'Mártir De Cuilapan'
$ws.Cells.Item(255, 2).Value = 'Taxco De Alarcón'
$ws.Cells.Item(257, 2).Value = 'Técpan De Galeana'
$ws.Cells.Item(263, 2).Value = 'Tlalixtaquilla De Maldonado'
$ws.Cells.Item(264, 2).Value = 'Tlapa De Comonfort'
$ws.Cells.Item(274, 2).Value = 'Agua Blanca De Iturbide'
$ws.Cells.Item(279, 2).Value = 'Atotonilco El Grande'
$ws.Cells.Item(282, 2).Value = 'Cuautepec De Hinojosa'
$ws.Cells.Item(291, 2).Value = 'Mixquiahuala De Juárez'
$ws.Cells.Item(292, 2).Value = 'Molango De Escamilla'
$ws.Cells.Item(293, 2).Value = 'Pachuca De Soto'
$ws.Cells.Item(295, 2).Value = 'Progreso De Obregón'
$ws.Cells.Item(302, 2).Value = 'Tula De Allende'
$ws.Cells.Item(303, 2).Value = 'Tulancingo De Bravo'
$ws.Cells.Item(304, 2).Value = 'Villa De Tezontepec'
$ws.Cells.Item(306, 2).Value = 'Zacualtipán De Ángeles'
$ws.Cells.Item(318, 2).Value = 'Encarnación De Díaz'
$ws.Cells.Item(322, 2).Value = 'Lagos De Moreno'
$ws.Cells.Item(324, 2).Value = 'Ojuelos De Jalisco'
$ws.Cells.Item(327, 2).Value = 'San Cristóbal De La Barranca'
$ws.Cells.Item(328, 2).Value = 'Talpa De Allende'
$ws.Cells.Item(329, 2).Value = 'Tamazula De Gordiano'
$ws.Cells.Item(335, 2).Value = 'Valle De Juárez'
$ws.Cells.Item(338, 2).Value = 'Zapotitlán De Vadillo'
$ws.Cells.Item(347, 2).Value = 'Coalcomán De Vázquez Pallares'
$ws.Cells.Item(379, 2).Value = 'Tiquicheo De Nicolás Romero'
$ws.Cells.Item(399, 2).Value = 'Jonacatepec De Leandro Valle'
$ws.Cells.Item(401, 2).Value = 'Puente De Ixtla'
$ws.Cells.Item(407, 2).Value = 'Tetela Del Volcán'
$ws.Cells.Item(409, 2).Value = 'Tlaltizapán De Zapata'
$ws.Cells.Item(416, 2).Value = 'Zacualpan De Amilpas'
$ws.Cells.Item(418, 2).Value = 'Bahía De Banderas'
$ws.Cells.Item(431, 2).Value = 'Acatlán De Pérez Figueroa'
$ws.Cells.Item(438, 2).Value = 'Capulálpam De Méndez'
$ws.Cells.Item(440, 2).Value = 'Chalcatongo De Hidalgo'
$ws.Cells.Item(441, 2).Value = 'Coicoyán De Las Flores'
$ws.Cells.Item(442, 2).Value = 'Constancia Del Rosario'
$ws.Cells.Item(445, 2).Value = 'Fresnillo De Trujano'
$ws.Cells.Item(446, 2).Value = 'Guadalupe De Ramírez'
$ws.Cells.Item(447, 2).Value = 'Guevea De Humboldt'
$ws.Cells.Item(448, 2).Value = 'Heroica Ciudad De Ejutla De Crespo'
$ws.Cells.Item(449, 2).Value = 'Heroica Ciudad De Huajuapan De León'
$ws.Cells.Item(450, 2).Value = 'Heroica Ciudad De Tlaxiaco'
$ws.Cells.Item(451, 2).Value = 'Huautla De Jiménez'
$ws.Cells.Item(452, 2).Value = 'Ixtlán De Juárez'
$ws.Cells.Item(453, 2).Value = 'Heroica Ciudad De Juchitán De Zaragoza'
$ws.Cells.Item(460, 2).Value = 'Mariscala De Juárez'
$ws.Cells.Item(461, 2).Value = 'Mártires De Tacubaya'
$ws.Cells.Item(463, 2).Value = 'Mazatlán Villa De Flores'
$ws.Cells.Item(465, 2).Value = 'Miahuatlán De Porfirio Díaz'
$ws.Cells.Item(467, 2).Value = 'Nejapa De Madero'
$ws.Cells.Item(469, 2).Value = 'Oaxaca De Juárez'
$ws.Cells.Item(470, 2).Value = 'Ocotlán De Morelos'
$ws.Cells.Item(471, 2).Value = 'Pinotepa De Don Luis'
$ws.Cells.Item(473, 2).Value = 'Putla Villa De Guerero'
$ws.Cells.Item(481, 2).Value = 'San Antonio De La Cal'
$ws.Cells.Item(489, 2).Value = 'San Felipe Jalapa De Díaz'
$ws.Cells.Item(495, 2).Value = 'San José Del Progreso'
$ws.Cells.Item(500, 2).Value = 'San Juan Bautista Lo De Soto'
$ws.Cells.Item(530, 2).Value = 'San Miguel Del Puerto'
$ws.Cells.Item(531, 2).Value = 'San Miguel Del Río'
$ws.Cells.Item(551, 2).Value = 'San Pedro El Alto'
$ws.Cells.Item(565, 2).Value = 'San Pedro Y San Pablo Teposcolula'
$ws.Cells.Item(575, 2).Value = 'Santa Ana Del Valle'
$ws.Cells.Item(591, 2).Value = 'Santa Cruz Tacache De Mina'
$ws.Cells.Item(599, 2).Value = 'Santa María Del Tule'
$ws.Cells.Item(604, 2).Value = 'Santa María Jalapa Del Marqués'
$ws.Cells.Item(634, 2).Value = 'Santo Domingo De Morelos'
$ws.Cells.Item(647, 2).Value = 'Tataltepec De Valdés'
$ws.Cells.Item(648, 2).Value = 'Teococuilco De Marcos Pérez'
$ws.Cells.Item(649, 2).Value = 'Teotitlán De Flores Magón'
$ws.Cells.Item(650, 2).Value = 'Tezoatlán De Segura Y Luna'
$ws.Cells.Item(651, 2).Value = 'Tlacolula De Matamoros'
$ws.Cells.Item(652, 2).Value = 'Tlalixtac De Cabrera'
$ws.Cells.Item(653, 2).Value = 'Totontepec Villa De Morelos'
$ws.Cells.Item(655, 2).Value = 'Villa De Etla'
$ws.Cells.Item(656, 2).Value = 'Villa De Tamazulápam Del Progreso'
$ws.Cells.Item(657, 2).Value = 'Villa De Tututepec De Melchor Ocampo'
$ws.Cells.Item(658, 2).Value = 'Villa De Zaachila'
$ws.Cells.Item(660, 2).Value = 'Villa Sola De Vega'
$ws.Cells.Item(661, 2).Value = 'Villa Talea De Castro'
$ws.Cells.Item(662, 2).Value = 'Zapotitlán Del Río'
$ws.Cells.Item(664, 2).Value = 'Zimatlán De Álvarez'
$ws.Cells.Item(687, 2).Value = 'Ayotoxco De Guerero'
$ws.Cells.Item(692, 2).Value = 'Chalchicomula De Sesma'
$ws.Cells.Item(703, 2).Value = 'Chila De La Sal'
$ws.Cells.Item(713, 2).Value = 'Cuayuca De Andrade'
$ws.Cells.Item(714, 2).Value = 'Cuetzalan Del Progreso'
$ws.Cells.Item(727, 2).Value = 'Huehuetlán El Chico'
$ws.Cells.Item(728, 2).Value = 'Huehuetlán El Grande'
$ws.Cells.Item(732, 2).Value = 'Ixcamilpa De Guerero'
$ws.Cells.Item(735, 2).Value = 'Izúcar De Matamoros'
$ws.Cells.Item(745, 2).Value = 'Los Reyes De Juárez'
$ws.Cells.Item(755, 2).Value = 'Palmar De Bravo'
$ws.Cells.Item(764, 2).Value = 'San Diego La Mesa Tochimiltzingo'
$ws.Cells.Item(781, 2).Value = 'San Nicolás De Los Ranchos'
$ws.Cells.Item(785, 2).Value = 'San Salvador El Seco'
$ws.Cells.Item(786, 2).Value = 'San Salvador El Verde'
$ws.Cells.Item(795, 2).Value = 'Tecali De Herrera'
$ws.Cells.Item(802, 2).Value = 'Tepanco De López'
$ws.Cells.Item(803, 2).Value = 'Tepatlaxco De Hidalgo'
$ws.Cells.Item(808, 2).Value = 'Tepexi De Rodríguez'
$ws.Cells.Item(810, 2).Value = 'Tetela De Ocampo'
$ws.Cells.Item(815, 2).Value = 'Tlacotepec De Benito Juárez'
$ws.Cells.Item(821, 2).Value = 'Totoltepec De Guerero'
$ws.Cells.Item(826, 2).Value = 'Xayacatlán De Bravo'
$ws.Cells.Item(832, 2).Value = 'Xochitlán De Vicente Suárez'
$ws.Cells.Item(844, 2).Value = 'Cadereyta De Montes'
$ws.Cells.Item(846, 2).Value = 'Jalpan De Serra'
$ws.Cells.Item(848, 2).Value = 'Pinal De Amoles'
$ws.Cells.Item(850, 2).Value = 'San Juan Del Río'
$ws.Cells.Item(856, 2).Value = 'Armadillo De Los Infante'
$ws.Cells.Item(857, 2).Value = 'Axtla De Terrazas'
$ws.Cells.Item(866, 2).Value = 'Villa De Arista'
$ws.Cells.Item(867, 2).Value = 'Villa De Reyes'
$ws.Cells.Item(884, 2).Value = 'Jalpa De Méndez'
$ws.Cells.Item(897, 2).Value = 'Acuamanala De Miguel Hidalgo'
$ws.Cells.Item(901, 2).Value = 'Contla De Juan Cuamatzi'
$ws.Cells.Item(906, 2).Value = 'Ixtacuixtla De Mariano Matamoros'
$ws.Cells.Item(908, 2).Value = 'Mazatecochco De José María Morelos'
$ws.Cells.Item(909, 2).Value = 'Muñoz De Domingo Arenas'
$ws.Cells.Item(912, 2).Value = 'Papalotla De Xicohténcatl'
$ws.Cells.Item(914, 2).Value = 'San Pablo Del Monte'
$ws.Cells.Item(936, 2).Value = 'Alto Lucero De Gutiérrez Barrios'
$ws.Cells.Item(939, 2).Value = 'Amatlán De Los Reyes'
$ws.Cells.Item(942, 2).Value = 'Boca Del Río'
$ws.Cells.Item(952, 2).Value = 'Cosamaloapan De Carpio'
$ws.Cells.Item(953, 2).Value = 'Cosautlán De Carvajal'
$ws.Cells.Item(965, 2).Value = 'Hueyapan De Ocampo'
$ws.Cells.Item(966, 2).Value = 'Huiloapan De Cuauhtémoc'
$ws.Cells.Item(967, 2).Value = 'Ignacio De La Llave'
$ws.Cells.Item(971, 2).Value = 'Ixhuatlán De Madero'
$ws.Cells.Item(972, 2).Value = 'Ixhuatlán Del Sureste'
$ws.Cells.Item(980, 2).Value = 'Juchique De Ferrer'
$ws.Cells.Item(984, 2).Value = 'Lerdo De Tejada'
$ws.Cells.Item(986, 2).Value = 'Martínez De La Torre'
$ws.Cells.Item(988, 2).Value = 'Medellín De Bravo'
$ws.Cells.Item(997, 2).Value = 'Paso De Ovejas'
$ws.Cells.Item(999, 2).Value = 'Poza Rica De Hidalgo'
$ws.Cells.Item(1028, 2).Value = 'Vega De Alatorre'
$ws.Cells.Item(1033, 2).Value = 'Zontecomatlán De López Y Fuentes'
$ws.Cells.Item(1034, 2).Value = 'Zozocolco De Hidalgo'
$ws.Cells.Item(1046, 2).Value = 'Tlaltenango De Sánchez Román'

# --- Remove the trailing footnote/source rows (1051-1055) ---
$ws.Range("A1051:A1055").EntireRow.Delete()

